# residential_electricity_price.xlsx — extend the monthly series through
# October 2024 (commit: "analaysis now till 2023").
#
# The sheet is a simple two-column (Date, Dollars) table that ran through
# row 265 (Dec-2022). We append 22 more monthly observations (Jan-2023 ..
# Oct-2024) in rows 266-287, copying the date/number formatting already
# used by the preceding data row so the new cells pick up the same styles
# (short-date column A, bordered "Dollars" column B) instead of generic
# formatting. Finally the sheet selection is left the way the author left
# it after entering the last value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = date serials (1st of each month), column B = price in $/kWh.
$dates  = @(44927,44958,44986,45017,45047,45078,45108,45139,45170,45200,45231,45261,45292,45323,45352,45383,45413,45444,45474,45505,45536,45566)
$prices = @(0.26469999999999999,0.27250000000000002,0.3039,0.29699999999999999,0.29820000000000002,0.31240000000000001,0.29459999999999997,0.29970000000000002,0.30009999999999998,0.31780000000000003,0.2954,0.29120000000000001,0.2949,0.31240000000000001,0.32479999999999998,0.34260000000000002,0.34300000000000003,0.32990000000000003,0.32550000000000001,0.3105,0.31640000000000001,0.30220000000000002)

$firstNewRow = 266

# Use the last existing data row as the formatting template.
$templateRow = $firstNewRow - 1
$ws.Range("A$templateRow`:B$templateRow").Copy()

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $firstNewRow + $i

    # Pulls in the date number format (col A) and the bordered $ style
    # (col B) from the template row without disturbing the values.
    $ws.Range("A$row`:B$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value = $dates[$i]
    $ws.Range("B$row").Value = $prices[$i]
}

$excel.CutCopyMode = 0

# Matches the author's final on-screen selection after the last edit.
$ws.Range("G283").Select()
